$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"
$meta.Range("B8").Value = "2024-06-19T17:47:42+02:00"
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("Z12").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAActClassObservation"
$elements.Range("Z13").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAActMood"

$elements.Range("F12").NumberFormat = "@"
$elements.Range("F12").Value = "0"
$elements.Range("F12").NumberFormat = "General"

$elements.Range("AG12").NumberFormat = "@"
$elements.Range("AG12").Value = "0"
$elements.Range("AG12").NumberFormat = "General"

$elements.Columns.Item(26).ColumnWidth = 58.0
